$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 through 15 (the extra dishes), keeping only the first
# dish ("Strogonoff do Chefe") and its 5 ingredient rows.
$ws.Range("A7:D15").EntireRow.Delete() | Out-Null

# Rename the dish from "[TESTE] Strogonoff do Chefe" to "Strogonoff do Chefe"
# for the remaining rows (2-6).
$ws.Range("A2:A6").Value = "Strogonoff do Chefe"

# Update the selected cell to A2 as in the saved file.
$ws.Range("A2").Select() | Out-Null
